# feat: add 2022-Q4 data
#
# - Duplicate the "2022-Q3" sheet, placing the copy right before it, and
#   rename the copy to "2022-Q4" (existing "2022-Q3"/"2022-Q2" sheets simply
#   shift right).
# - Overwrite the values on the new "2022-Q4" sheet with the Q4 numbers.
# - Insert a new row 2 on the "总计" (total) sheet summarizing 2022-Q4 and
#   push the existing 2022-Q3 / 2022-Q2 summary rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate "2022-Q3" -> "2022-Q4" (inserted before "2022-Q3")
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2) Update the data on the new "2022-Q4" sheet
# ---------------------------------------------------------------------
$q4.Range("B2").Value = "'011550"
$q4.Range("C2").Value = "湘财创新成长一年持有期混合A"
$q4.Range("D2").Value = "'1.23"
$q4.Range("E2").Value = "'91.30"
$q4.Range("F2").Value = "'4.70"
$q4.Range("G2").Value = "'0.0578"
$q4.Range("H2").Value = 2

$q4.Range("D3").Value = "'1.70"
$q4.Range("E3").Value = "'89.45"
$q4.Range("F3").Value = "'3.39"
$q4.Range("G3").Value = "'0.0576"
$q4.Range("H3").Value = 3

$q4.Range("B4").Value = "'005075"
$q4.Range("C4").Value = "富国研究量化精选混合"
$q4.Range("D4").Value = "'2.48"
$q4.Range("E4").Value = "'90.71"
$q4.Range("F4").Value = "'1.85"
$q4.Range("G4").Value = "'0.0459"
$q4.Range("H4").Value = 1

$q4.Range("B5").Value = "'011551"
$q4.Range("C5").Value = "湘财创新成长一年持有期混合C"
$q4.Range("D5").Value = "'0.16"
$q4.Range("E5").Value = "'91.30"
$q4.Range("F5").Value = "'4.70"
$q4.Range("G5").Value = "'0.0075"
$q4.Range("H5").Value = 2

$q4.Range("D6").Value = "'0.15"
$q4.Range("E6").Value = "'89.45"
$q4.Range("F6").Value = "'3.39"
$q4.Range("G6").Value = "'0.0051"
$q4.Range("H6").Value = 3

# ---------------------------------------------------------------------
# 3) Update the "总计" (total) sheet: insert a 2022-Q4 row before the
#    existing 2022-Q3 row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.17

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
